# Auto-generated Excel COM-interop script
# Applies numeric cell updates (value changes, clears, and additions)
# derived from the Masamune_Profits.xlsx diff, across 7 worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H75").Value = 42782.5
$ws.Range("J75").Value = 46037.145
$ws.Range("L75").Value = 46037.145
$ws.Range("N75").Value = -47909.145
$ws.Range("H78").Value = 42782.5
$ws.Range("J78").Value = 46037.145
$ws.Range("L78").Value = 138111.435
$ws.Range("N78").Value = -147471.435
$ws.Range("H108").Value = 32148.5
$ws.Range("J108").Value = 32148.5
$ws.Range("L108").Value = 32148.5
$ws.Range("N108").Value = -39828.5
$ws.Range("H120").Value = 49000
$ws.Range("J120").Value = 49000
$ws.Range("L120").Value = 49000
$ws.Range("N120").Value = -58676
$ws.Range("H126").Value = 55019.715
$ws.Range("J126").Value = 55019.715
$ws.Range("L126").Value = 55019.715
$ws.Range("N126").Value = -64899.715
$ws.Range("H134").Value = 54950
$ws.Range("J134").Value = 54950
$ws.Range("L134").Value = 54950
$ws.Range("N134").Value = -65090
$ws.Range("H138").Value = 3718.074
$ws.Range("I138").Value = 3084.5789
$ws.Range("J138").Value = 3912.2097
$ws.Range("K138").Value = 9253.736699999999
$ws.Range("L138").Value = 11736.6291
$ws.Range("M138").Value = -4113.736699999999
$ws.Range("N138").Value = -22016.6291

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H80").Value = 62495
$ws.Range("J80").Value = 62495
$ws.Range("L80").Value = 62495
$ws.Range("N80").Value = -64491
$ws.Range("H83").Value = 62495
$ws.Range("J83").Value = 62495
$ws.Range("L83").Value = 187485
$ws.Range("N83").Value = -197469
$ws.Range("H109").Value = 27753.666
$ws.Range("J109").Value = 27753.666
$ws.Range("L109").Value = 27753.666
$ws.Range("N109").Value = -30527.666
$ws.Range("H118").Value = 49997.332
$ws.Range("J118").Value = 49997.332
$ws.Range("L118").Value = 49997.332
$ws.Range("N118").Value = -53311.332
$ws.Range("H120").Value = 47366
$ws.Range("J120").Value = 47366
$ws.Range("L120").Value = 47366
$ws.Range("N120").Value = -57042
$ws.Range("H122").Value = 1507.6923
$ws.Range("I122").Value = 1533.3334
$ws.Range("J122").Value = 1450
$ws.Range("K122").Value = 4600.0002
$ws.Range("L122").Value = 4350
$ws.Range("M122").Value = -2150.0002
$ws.Range("N122").Value = -9250
$ws.Range("H128").Value = 46947.332
$ws.Range("J128").Value = 46947.332
$ws.Range("L128").Value = 46947.332
$ws.Range("N128").Value = -56907.332

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H70").Value = 200000
$ws.Range("J70").Value = 200000
$ws.Range("L70").Value = 200000
$ws.Range("N70").Value = -200586
$ws.Range("H73").Value = 200000
$ws.Range("J73").Value = 200000
$ws.Range("L73").Value = 200000
$ws.Range("N73").Value = -202028
$ws.Range("H126").Value = 50936.25
$ws.Range("J126").Value = 50936.25
$ws.Range("L126").Value = 50936.25
$ws.Range("N126").Value = -60816.25
$ws.Range("H130").Value = 41336
$ws.Range("J130").Value = 41336
$ws.Range("L130").Value = 41336
$ws.Range("N130").Value = -51376

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H20").Value = 47847.668
$ws.Range("J20").Value = 47847.668
$ws.Range("L20").Value = 47847.668
$ws.Range("N20").Value = -48319.668
$ws.Range("H30").Value = 47847.668
$ws.Range("J30").Value = 47847.668
$ws.Range("L30").Value = 47847.668
$ws.Range("N30").Value = -48029.668
$ws.Range("H31").Value = 4991.4473
$ws.Range("I31").Value = 2175.423
$ws.Range("J31").Value = 11092.833
$ws.Range("K31").Value = 2175.423
$ws.Range("L31").Value = 11092.833
$ws.Range("M31").Value = -1880.423
$ws.Range("N31").Value = -11682.833
$ws.Range("H34").Value = 4991.4473
$ws.Range("I34").Value = 2175.423
$ws.Range("J34").Value = 11092.833
$ws.Range("K34").Value = 2175.423
$ws.Range("L34").Value = 11092.833
$ws.Range("M34").Value = -1973.423
$ws.Range("N34").Value = -11496.833
$ws.Range("H70").Value = 34851
$ws.Range("J70").Value = 34851
$ws.Range("L70").Value = 34851
$ws.Range("N70").Value = -35481
$ws.Range("H73").Value = 34851
$ws.Range("J73").Value = 34851
$ws.Range("L73").Value = 34851
$ws.Range("N73").Value = -37035
$ws.Range("H82").Value = 46200
$ws.Range("J82").Value = 46200
$ws.Range("L82").Value = 46200
$ws.Range("N82").Value = -46922
$ws.Range("H85").Value = 46200
$ws.Range("J85").Value = 46200
$ws.Range("L85").Value = 46200
$ws.Range("N85").Value = -48696
$ws.Range("H97").Value = 19381
$ws.Range("J97").Value = 19381
$ws.Range("L97").Value = 19381
$ws.Range("N97").Value = -21363
$ws.Range("H116").Value = 64991
$ws.Range("J116").Value = 64991
$ws.Range("L116").Value = 64991
$ws.Range("N116").Value = -74169
$ws.Range("H118").Value = 44742
$ws.Range("J118").Value = 44742
$ws.Range("L118").Value = 44742
$ws.Range("N118").Value = -48056
$ws.Range("H128").Value = 47847.668
$ws.Range("J128").Value = 47847.668
$ws.Range("L128").Value = 47847.668
$ws.Range("N128").Value = -57807.668
$ws.Range("H135").Value = 57437
$ws.Range("J135").Value = 57437
$ws.Range("L135").Value = 57437
$ws.Range("N135").Value = -67577

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 46639
$ws.Range("J98").Value = 46639
$ws.Range("L98").Value = 46639
$ws.Range("N98").Value = -52629
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").ClearContents()
$ws.Range("N101").Value = 0
$ws.Range("H130").Value = 44751.2
$ws.Range("J130").Value = 44751.2
$ws.Range("L130").Value = 44751.2
$ws.Range("N130").Value = -54791.2
$ws.Range("H136").Value = 36081.5
$ws.Range("J136").Value = 36081.5
$ws.Range("L136").Value = 108244.5
$ws.Range("N136").Value = -113344.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").ClearContents()
$ws.Range("N63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").ClearContents()
$ws.Range("N66").Value = 0
$ws.Range("H74").Value = 30217
$ws.Range("I74").Value = 0
$ws.Range("J74").Value = 30217
$ws.Range("K74").Value = 0
$ws.Range("L74").ClearContents()
$ws.Range("M74").Value = 30217
$ws.Range("N74").Value = -32213
$ws.Range("H77").Value = 30217
$ws.Range("I77").Value = 0
$ws.Range("J77").Value = 30217
$ws.Range("K77").Value = 0
$ws.Range("L77").ClearContents()
$ws.Range("M77").Value = 90651
$ws.Range("N77").Value = -100635
$ws.Range("H81").Value = 20000
$ws.Range("I81").Value = 0
$ws.Range("K81").Value = 0
$ws.Range("M81").ClearContents()
$ws.Range("H84").Value = 20000
$ws.Range("I84").Value = 0
$ws.Range("K84").Value = 0
$ws.Range("M84").ClearContents()
$ws.Range("H88").Value = 31120.666
$ws.Range("H91").Value = 31120.666
$ws.Range("H92").Value = 24596
$ws.Range("J92").Value = 24596
$ws.Range("L92").Value = 24596
$ws.Range("N92").Value = -29588
$ws.Range("H99").Value = 25500
$ws.Range("I99").Value = 16500
$ws.Range("K99").Value = 16500
$ws.Range("M99").Value = -13505
$ws.Range("H102").Value = 29800
$ws.Range("J102").Value = 29800
$ws.Range("L102").Value = 29800
$ws.Range("N102").Value = -36290
$ws.Range("H111").Value = 44383
$ws.Range("J111").Value = 44383
$ws.Range("L111").Value = 44383
$ws.Range("N111").Value = -52563
$ws.Range("H127").Value = 50715
$ws.Range("J127").Value = 50715
$ws.Range("L127").Value = 50715
$ws.Range("N127").Value = -60635
$ws.Range("H130").Value = 44141.6
$ws.Range("J130").Value = 44141.6
$ws.Range("L130").Value = 44141.6
$ws.Range("N130").Value = -54181.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 44210
$ws.Range("J16").Value = 44210
$ws.Range("L16").Value = 44210
$ws.Range("N16").Value = -44794
$ws.Range("H75").Value = 22318.5
$ws.Range("J75").Value = 22318.5
$ws.Range("L75").Value = 22318.5
$ws.Range("N75").Value = -24190.5
$ws.Range("H78").Value = 22318.5
$ws.Range("J78").Value = 22318.5
$ws.Range("L78").Value = 66955.5
$ws.Range("N78").Value = -76315.5
$ws.Range("H93").Value = 33108.332
$ws.Range("J93").Value = 33108.332
$ws.Range("L93").Value = 33108.332
$ws.Range("N93").Value = -38100.332
$ws.Range("H119").Value = 47897.332
$ws.Range("J119").Value = 47897.332
$ws.Range("L119").Value = 47897.332
$ws.Range("N119").Value = -57573.332
$ws.Range("H123").Value = 28406.154
$ws.Range("J123").Value = 28406.154
$ws.Range("L123").Value = 28406.154
$ws.Range("N123").Value = -38206.15399999999
$ws.Range("H135").Value = 49404.668
$ws.Range("J135").Value = 49404.668
$ws.Range("L135").Value = 49404.668
$ws.Range("N135").Value = -59544.668
